# Update gh-pages output numbers (view/attendee counts) and a date-range
# string, as generated at commit 456a3b4.
#
# Sheet "展览" (Exhibition) -- F column view counts.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 1270
$ws1.Range("F8").Value  = 7564
$ws1.Range("F12").Value = 8207
$ws1.Range("F16").Value = 5616
$ws1.Range("F17").Value = 5616
$ws1.Range("F19").Value = 2568
$ws1.Range("F26").Value = 484
$ws1.Range("F27").Value = 2790
$ws1.Range("F28").Value = 2790
$ws1.Range("F30").Value = 4
$ws1.Range("F31").Value = 2809
$ws1.Range("F32").Value = 2809
$ws1.Range("F36").Value = 276
$ws1.Range("F39").Value = 866
$ws1.Range("F40").Value = 1625
$ws1.Range("F43").Value = 2615
$ws1.Range("F45").Value = 2271
$ws1.Range("F46").Value = 7

# Sheet "演出" (Performance) -- row 8 date range extended + attendee count.
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("E8").Value = "2024.05.25 19:30-06.06 22:00"
$ws2.Range("G8").Value = 480

# Sheet "本地生活" (Local life) -- F column view count.
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1308

# Sheet "全部类型" (All types) -- aggregated view, mirrors the other sheets.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1308
$ws4.Range("F6").Value  = 1270
$ws4.Range("F7").Value  = 7564
$ws4.Range("F11").Value = 8207
$ws4.Range("F15").Value = 5616
$ws4.Range("F16").Value = 5616
$ws4.Range("F18").Value = 2568
$ws4.Range("F26").Value = 484
$ws4.Range("F27").Value = 2790
$ws4.Range("F28").Value = 2790
$ws4.Range("F30").Value = 4
$ws4.Range("F31").Value = 2809
$ws4.Range("F32").Value = 2809
$ws4.Range("F35").Value = 276
$ws4.Range("F40").Value = 866
$ws4.Range("F42").Value = 1625
$ws4.Range("F45").Value = 2615
$ws4.Range("F48").Value = 2271
$ws4.Range("F49").Value = 7
$ws4.Range("E52").Value = "2024.05.25 19:30-06.06 22:00"
$ws4.Range("G52").Value = 480
